$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds price strings that look like numbers (e.g. "29,499"). Force
# that range to Text format first so Excel keeps them as literal text
# instead of silently converting them to numbers.
$ws.Range("A2:AC2").NumberFormat = "@"

# Row 1: product names
$ws.Cells.Item(1,1).Value = "Apple iPhone 7 (32GB) - Black"
$ws.Cells.Item(1,2).Value = "Apple iPhone 11 (128GB) - Black"
$ws.Cells.Item(1,3).Value = "Vivo U10 (Electric Blue, 5000 mAH 18W Fast Charge Battery, 3GB RAM, 32GB Storage)"
$ws.Cells.Item(1,4).Value = "Nokia 105 2019 (Single SIM, Black)"
$ws.Cells.Item(1,5).Value = "Samsung Galaxy M21 (Midnight Blue, 4GB RAM, 64GB Storage)"
$ws.Cells.Item(1,6).Value = "Samsung Galaxy M31 (Ocean Blue, 6GB RAM, 128GB Storage)"
$ws.Cells.Item(1,7).Value = "Vivo Y91i (Fusion Black, 2GB RAM, 32GB Storage) with No Cost EMI/Additional Exchange Offers"
$ws.Cells.Item(1,8).Value = "Redmi Note 8 Pro (Halo White, 6GB RAM, 128GB Storage with Helio G90T Processor)"
$ws.Cells.Item(1,9).Value = "SNEXIAN Guru 106 Dual Sim Basic Mobile Phone with Digital Camera and 1.8 inch Screen (Black, Upto 16GB) - NO Charger"
$ws.Cells.Item(1,10).Value = "Vivo Y91i (Ocean Blue, 2GB RAM, 32GB Storage) with No Cost EMI/Additional Exchange Offers"
$ws.Cells.Item(1,11).Value = "Karbonn KX3 1.8 inch Display Feature Phone with Bluetooth,Dual Sim, 0.3 MP Digital Camera with Zoom, 800 mAH Battery,32 GB Expandable Memory and Support for MP3+MP4, Boom Box Speaker,Black Red Colour."
$ws.Cells.Item(1,12).Value = "itel A46 (Neon Water, 2GB RAM, 16GB Storage)"
$ws.Cells.Item(1,13).Value = "Samsung Galaxy M21 (Raven Black, 6GB RAM, 128GB Storage)"
$ws.Cells.Item(1,14).Value = "SNEXIAN Guru 311 Dual Sim Basic Mobile Phone with Digital Camera and 1.8 inch Screen (Black, Upto 16GB) - NO Charger"
$ws.Cells.Item(1,15).Value = "Maono AU-D30 BassCurve Neck Band in-Ear Bluetooth Wireless Earphones, with Bluetooth 5.0, Sweatproof Headphones, Long Life Battery, Flexible Headset and Built-in Mic"
$ws.Cells.Item(1,16).Value = "Spigen [Rugged Armor] Case for Huawei P30, Patent Design Flexible TPU Phone Case Cover for Huawei P30 (NOT Compatible with P30 Pro / P30 Lite)"
$ws.Cells.Item(1,17).Value = "TECNO Spark 5 Pro (Seabed Blue, 4GB RAM, 64GB Storage)"
$ws.Cells.Item(1,18).Value = "Nokia 105 2019 (Single SIM, Blue)"
$ws.Cells.Item(1,19).Value = "Nokia 105 2019 (Single SIM, Pink)"
$ws.Cells.Item(1,20).Value = "ELV 4mm Thickness Aluminum mobile Stand (3.5 - 8 inches) - Black"
$ws.Cells.Item(1,21).Value = "IKall K76 Mobile (1.4 Inch Display, Single Sim, 600 mAh Battery) (Yellow)"
$ws.Cells.Item(1,22).Value = "Samsung Galaxy M21 (Raven Black, 4GB RAM, 64GB Storage)"
$ws.Cells.Item(1,23).Value = "OPPO A5 2020 (Dazzling White, 3GB RAM, 64GB Storage) with No Cost EMI/Additional Exchange Offers"
$ws.Cells.Item(1,24).Value = "OPPO F11 (Fluorite Purple, 6GB RAM, 128GB Storage) with No Cost EMI/Additional Exchange Offers"
$ws.Cells.Item(1,25).Value = "OPPO A5 2020 (Dazzling White, 4GB RAM, 64GB Storage) with No Cost EMI/Additional Exchange Offers"
$ws.Cells.Item(1,26).Value = "Vivo U10 (Thunder Black,5000 mAH 18W Fast Charge Battery, 3GB RAM, 32GB Storage)"
$ws.Cells.Item(1,27).Value = "Micromax X421 Black"
$ws.Cells.Item(1,28).Value = "Redmi 8A Dual (Sky White, 2GB RAM, 32GB Storage) – Dual Cameras & 5,000 mAH Battery"
$ws.Cells.Item(1,29).Value = "Lenovo HT10 True Wireless Earbuds Earphones Headphones (Bluetooth V5.0) in-Built Mic with Extra HD Sound AirBass Rated IPX5 Waterproof and Sweatproof (Black)"
$ws.Cells.Item(1,30).Value = "Spigen Slim Armor CS Flip Wallet Card Slot Holder Back Cover Case Designed for iPhone 11 - Gunmetal"

# Row 2: prices (kept as text, matching the source formatting with commas)
$ws.Cells.Item(2,1).Value = "29,499"
$ws.Cells.Item(2,2).Value = "73,600"
$ws.Cells.Item(2,3).Value = "10,990"
$ws.Cells.Item(2,4).Value = "1,241"
$ws.Cells.Item(2,5).Value = "13,999"
$ws.Cells.Item(2,6).Value = "17,499"
$ws.Cells.Item(2,7).Value = "7,990"
$ws.Cells.Item(2,8).Value = "16,999"
$ws.Cells.Item(2,9).Value = "599"
$ws.Cells.Item(2,10).Value = "7,990"
$ws.Cells.Item(2,11).Value = "805"
$ws.Cells.Item(2,12).Value = "15,999"
$ws.Cells.Item(2,13).Value = "599"
$ws.Cells.Item(2,14).Value = "1,299"
$ws.Cells.Item(2,15).Value = "899"
$ws.Cells.Item(2,16).Value = "10,499"
$ws.Cells.Item(2,17).Value = "1,105"
$ws.Cells.Item(2,18).Value = "1,230"
$ws.Cells.Item(2,19).Value = "149"
$ws.Cells.Item(2,20).Value = "399"
$ws.Cells.Item(2,21).Value = "13,999"
$ws.Cells.Item(2,22).Value = "10,990"
$ws.Cells.Item(2,23).Value = "18,990"
$ws.Cells.Item(2,24).Value = "11,990"
$ws.Cells.Item(2,25).Value = "10,990"
$ws.Cells.Item(2,26).Value = "1,020"
$ws.Cells.Item(2,27).Value = "7,499"
$ws.Cells.Item(2,28).Value = "3,999"
$ws.Cells.Item(2,29).Value = "1,599"
